$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: set column C (Ordinal flag) to 0 for all data rows (2-82)
$ws.Range("C2:C82").Value = 0

# Step 2: set column C to 1 for the ordinal-feature rows
$ordinalRows = @(19,20,21,22,29,30,32,33,42,55,59,61,65,66,74)
foreach ($r in $ordinalRows) {
    $ws.Cells.Item($r, 3).Value = 1
}

# Step 3: set column D notes, in row order, so shared strings are created
#         in the same order as the original authoring session
#         (Quality has order, Condition has order, Year has order, then Remove Row)
$ws.Cells.Item(19, 4).Value = "Quality has order"
$ws.Cells.Item(20, 4).Value = "Condition has order"
$ws.Cells.Item(21, 4).Value = "Year has order"
$ws.Cells.Item(22, 4).Value = "Year has order"
$ws.Cells.Item(29, 4).Value = "Quality has order"
$ws.Cells.Item(30, 4).Value = "Condition has order"
$ws.Cells.Item(32, 4).Value = "Quality has order"
$ws.Cells.Item(33, 4).Value = "Condition has order"
$ws.Cells.Item(42, 4).Value = "Quality has order"
$ws.Cells.Item(55, 4).Value = "Quality has order"
$ws.Cells.Item(59, 4).Value = "Quality has order"
$ws.Cells.Item(61, 4).Value = "Year has order"
$ws.Cells.Item(65, 4).Value = "Quality has order"
$ws.Cells.Item(66, 4).Value = "Condition has order"
$ws.Cells.Item(74, 4).Value = "Quality has order"

# Step 4: the "Remove Row" note on row 2 is added last so it becomes the last shared string
$ws.Cells.Item(2, 4).Value = "Remove Row"

# Step 5: size the new column D and select D15 like the saved workbook
$ws.Columns.Item(4).ColumnWidth = 16.7
$ws.Range("D15").Select()
